# Improve inc implementation, delete unnecessary files.
# Removes the rows for retired/outdated incentive entries (the MI
# "Alternative Fuel Development Property Tax Exemption", the AL utility
# "Biofuel Production Jobs Tax Credit" duplicate, and the two US "Second
# Generation Biofuel ..." entries) from the bottom of the incentives table,
# leaving the now-unused placeholder cells blank (same as the surviving
# rows around them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 (was "Alternative Fuel Development Property Tax Exemption", MI) ---
# D23 keeps its style (an exemption-only formatted cell) but the value goes away.
$ws.Range("A23").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Range("D23").ClearContents()
$ws.Range("E23").Clear()
$ws.Range("F23").Clear()
$ws.Range("H23").Clear()

# --- Row 24 (was "Biofuel Production Jobs Tax Credit", AL / utility) ---
# Entire row's data cells are removed outright.
$ws.Range("A24").Clear()
$ws.Range("B24").Clear()
$ws.Range("C24").Clear()
$ws.Range("D24").Clear()
$ws.Range("E24").Clear()
$ws.Range("F24").Clear()

# --- Row 25 (was "Second Generation Biofuel Producer Tax Credit", US) ---
$ws.Range("A25").Clear()
$ws.Range("B25").Clear()
$ws.Range("C25").Clear()
$ws.Range("D25").ClearContents()
$ws.Range("E25").Clear()
$ws.Range("F25").Clear()
$ws.Range("H25").Clear()
$ws.Range("I25").Clear()

# --- Row 26 (was "Second Generation Biofuel Plant Depreciation Deduction
#     Allowance", US) ---
$ws.Range("A26").Clear()
$ws.Range("B26").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("I26").Clear()

# Update the view's selection to the now-edited block.
$ws.Range("A22:I26").Select()
